$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-02-24 Saturday"; new = "2024-02-25 Sunday"},
    @{old = "442÷5="; new = "387÷3="},
    @{old = "121÷7="; new = "602÷3="},
    @{old = "886÷8="; new = "754÷5="},
    @{old = "278÷5="; new = "729÷2="},
    @{old = "267÷4="; new = "352÷6="},
    @{old = "486÷6="; new = "763÷6="},
    @{old = "748÷9="; new = "397÷7="},
    @{old = "875÷7="; new = "251÷7="},
    @{old = "764÷7="; new = "504÷4="},
    @{old = "701÷6="; new = "923÷3="},
    @{old = "894÷8="; new = "808÷5="},
    @{old = "950÷8="; new = "865÷5="},
    @{old = "607÷6="; new = "119÷6="},
    @{old = "998÷8="; new = "487÷6="},
    @{old = "147÷5="; new = "225÷3="},
    @{old = "375÷2="; new = "970÷9="},
    @{old = "377÷8="; new = "930÷8="},
    @{old = "242÷3="; new = "706÷2="},
    @{old = "779÷7="; new = "550÷2="},
    @{old = "139÷6="; new = "136÷4="},
    @{old = "426÷5="; new = "598÷6="},
    @{old = "349÷3="; new = "972÷5="},
    @{old = "790÷2="; new = "417÷8="},
    @{old = "350÷9="; new = "872÷6="},
    @{old = "995÷4="; new = "511÷7="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
